# Alteração nos rótulos da tabela para já transformar a primeira linha em
# cabeçalho automaticamente no Power BI.
#
# Sheets 1, 2, 3, 5 and 6 use a "year" header row (2015, 2030, 2040, 2050),
# which gets prefixed with "Ano ".
# Sheet 4 uses an "interval" header row (2015, 2015-2030, 2031-2040,
# 2041-2050), which gets prefixed with "Intervalo ".

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)",
    "Custo Total (bilhões de R$)"
)

foreach ($name in $anoSheets) {
    $ws = $wb.Worksheets.Item($name)
    $lastCol = $ws.Cells.Item(1, 1).End(-4161).Column   # xlToRight = -4161
    for ($col = 2; $col -le $lastCol; $col++) {
        $cell = $ws.Cells.Item(1, $col)
        $txt = $cell.Text.ToString()
        if (-not $txt.StartsWith("Ano ")) {
            $cell.Value = "Ano " + $txt
        }
    }
}

$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$lastCol = $wsIntervalo.Cells.Item(1, 1).End(-4161).Column
for ($col = 2; $col -le $lastCol; $col++) {
    $cell = $wsIntervalo.Cells.Item(1, $col)
    $txt = $cell.Text.ToString()
    if (-not $txt.StartsWith("Intervalo ")) {
        $cell.Value = "Intervalo " + $txt
    }
}
